# [base commands] - [assertMatch(text,regex)]: NEW command to check for text
# value via regular expression.
#
# The workbook's hidden "#system" sheet stores, per category-column, the
# catalog of available step commands; column A ("target") lists the category
# names themselves. This script:
#   1) Inserts the new "assertMatch(text,regex)" command into the "base"
#      category (column F), in alphabetical order.
#   2) Inserts the new "openFile(filePath)" command into the "external"
#      category (column J), in alphabetical order.
#   3) Removes the "tn.5250" entry from the "target" category list
#      (column A), shifting subsequent category names up by one row.
#   4) Shifts the "web"/"webalert"/"webcookie"/"ws"/"ws.async"/"xml" category
#      data left by one column (AB->AA, AC->AB, AD->AC, AE->AD, AF->AE,
#      AG->AF), clearing the vacated AG column.
#   5) Updates the affected named ranges accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) column F ("base"): insert "assertMatch(text,regex)" at F11, pushing
#    the existing F11:F44 block down to F12:F45.
# ---------------------------------------------------------------------
for ($r = 44; $r -ge 11; $r--) {
    $v = $ws.Range("F" + $r).Value2
    $ws.Range("F" + ($r + 1)).Value = $v
}
$ws.Range("F11").Value = "assertMatch(text,regex)"

# ---------------------------------------------------------------------
# 2) column J ("external"): insert "openFile(filePath)" at J2, pushing
#    the existing J2:J6 block down to J3:J7.
# ---------------------------------------------------------------------
for ($r = 6; $r -ge 2; $r--) {
    $v = $ws.Range("J" + $r).Value2
    $ws.Range("J" + ($r + 1)).Value = $v
}
$ws.Range("J2").Value = "openFile(filePath)"

# ---------------------------------------------------------------------
# 3) column A ("target"): drop the "tn.5250" row (A27), pulling A28:A33
#    up to A27:A32 and clearing the now-vacant A33.
# ---------------------------------------------------------------------
for ($r = 28; $r -le 33; $r++) {
    $v = $ws.Range("A" + $r).Value2
    $ws.Range("A" + ($r - 1)).Value = $v
}
$ws.Range("A33").ClearContents()

# ---------------------------------------------------------------------
# 4) columns AA..AG: shift the "web".."xml" catalogs one column to the
#    left (AA<-AB, AB<-AC, AC<-AD, AD<-AE, AE<-AF, AF<-AG), then clear AG.
#    Processed left-to-right so every source column is read before it is
#    itself overwritten.
# ---------------------------------------------------------------------
$destCols = @("AA", "AB", "AC", "AD", "AE", "AF")
$srcCols  = @("AB", "AC", "AD", "AE", "AF", "AG")

for ($i = 0; $i -lt $destCols.Length; $i++) {
    $destCol = $destCols[$i]
    $srcCol = $srcCols[$i]
    for ($r = 1; $r -le 151; $r++) {
        $v = $ws.Range($srcCol + $r).Value2
        if ($v -eq $null) {
            $ws.Range($destCol + $r).ClearContents()
        } else {
            $ws.Range($destCol + $r).Value = $v
        }
    }
}
for ($r = 1; $r -le 151; $r++) {
    $ws.Range("AG" + $r).ClearContents()
}

# ---------------------------------------------------------------------
# 5) Update named ranges to reflect the new row/column boundaries.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"
